$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1048539867316382
$ws.Range("H2").Value = 25.02302909863244
$ws.Range("I2").Value = -73.97306565948855
$ws.Range("G3").Value = 0.09101976076376014
$ws.Range("H3").Value = -21.77457158216032
$ws.Range("G4").Value = -0.03485151532418401
$ws.Range("H4").Value = -125.9376692778586
$ws.Range("G5").Value = -0.06184762211758178
$ws.Range("H5").Value = 13.92319949202543
$ws.Range("G6").Value = 0.01552484453078491
$ws.Range("H6").Value = -56.32422545406101
$ws.Range("G7").Value = 0.001583474222683695
$ws.Range("H7").Value = -92.19775042377823
$ws.Range("G8").Value = -0.1402592539175551
$ws.Range("H8").Value = 0.5943694742114293
$ws.Range("G9").Value = -0.1271473925996481
$ws.Range("H9").Value = 7.11507805088781
$ws.Range("G10").Value = -0.03246178750902648
$ws.Range("H10").Value = 69.95033602096809
$ws.Range("G11").Value = -0.0708121977820121
$ws.Range("H11").Value = -6.419117576433695
$ws.Range("G12").Value = -0.3157292477428361
$ws.Range("H12").Value = 23.95663691185233
$ws.Range("G13").Value = -0.467953977649264
$ws.Range("H13").Value = -4.24473503890964
$ws.Range("G14").Value = -0.08385027641122184
$ws.Range("H14").Value = -65.26880086519429
$ws.Range("G15").Value = 0.02058044575925375
$ws.Range("H15").Value = 124.8882428512505
$ws.Range("G16").Value = 0.11911472718085
$ws.Range("H16").Value = -18.42724537486222
$ws.Range("G17").Value = 0.1474446849492863
$ws.Range("H17").Value = 20.24376038809976
$ws.Range("G18").Value = 0.1372178443256943
$ws.Range("H18").Value = -0.8951044018356646
$ws.Range("G19").Value = 0.097961454324279
$ws.Range("H19").Value = 2.67770325920964
$ws.Range("G20").Value = 0.02764477448937599
$ws.Range("H20").Value = 7.843950436072178
$ws.Range("G21").Value = 0.04241477322348758
$ws.Range("H21").Value = -43.57269848764449
$ws.Range("G24").Value = 0.09187565468843437
$ws.Range("H24").Value = -8.523645974971302
$ws.Range("G25").Value = 0.2011942122054161
$ws.Range("H25").Value = 32.76068207599671
$ws.Range("G26").Value = 0.07060795724642746
$ws.Range("H26").Value = -10.76370969592137
$ws.Range("G27").Value = 0.05825412313763535
$ws.Range("H27").Value = -41.65989584513596
$ws.Range("G28").Value = -0.2324902651729035
$ws.Range("H28").Value = -9.089318090416343
$ws.Range("G29").Value = -0.201775765398231
$ws.Range("H29").Value = 1.71982252160831
$ws.Range("G30").Value = 0.05660815975405609
$ws.Range("H30").Value = 28.27113686696355
$ws.Range("G31").Value = 0.02997592756123089
$ws.Range("H31").Value = 13.83359372385789
$ws.Range("G32").Value = 0.08776827616742577
$ws.Range("H32").Value = -7.573919611653512
$ws.Range("G33").Value = 0.1107265907508278
$ws.Range("H33").Value = 6.503834124361433
$ws.Range("G34").Value = 0.04539576030892353
$ws.Range("H34").Value = -2.222662734409952
$ws.Range("G35").Value = 0.03150228623555582
$ws.Range("H35").Value = 315.7879193247942
$ws.Range("G36").Value = 0.07279564033479483
$ws.Range("H36").Value = 26.08308552718183
$ws.Range("G37").Value = 0.06352862012912477
$ws.Range("H37").Value = -9.665733636216741
$ws.Range("G38").Value = 0.02316099129221889
$ws.Range("H38").Value = -55.78816970431294
$ws.Range("G39").Value = 0.006397086235182236
$ws.Range("H39").Value = -69.15183177842735
$ws.Range("G40").Value = 0.0005753516362875627
$ws.Range("H40").Value = 106.7760257491479
$ws.Range("G41").Value = 0.009819843068733552
$ws.Range("H41").Value = -72.22569403736594
$ws.Range("G42").Value = 0.1322269571980616
$ws.Range("H42").Value = -1.091488423737454
$ws.Range("G43").Value = 0.1449752262303823
$ws.Range("H43").Value = -2.685734129487327
$ws.Range("G44").Value = 0.009475087163798891
$ws.Range("H44").Value = 211.3212478751
$ws.Range("G45").Value = 0.005084800236594285
$ws.Range("H45").Value = 146.3102252068789
$ws.Range("G46").Value = -0.01057402027488038
$ws.Range("H46").Value = -221.1015141608299
$ws.Range("G47").Value = -0.02601085296065285
$ws.Range("H47").Value = -180.3352476946778
$ws.Range("G48").Value = 0.06111628288279508
$ws.Range("H48").Value = 21.57094746220717
$ws.Range("G49").Value = 0.07175347568865321
$ws.Range("H49").Value = 8.610778303351996
$ws.Range("G50").Value = 0.1364161820132598
$ws.Range("H50").Value = -15.40484661302834
$ws.Range("G51").Value = 0.1529660957478109
$ws.Range("H51").Value = -10.60624026090028
$ws.Range("G52").Value = -0.1880782984102244
$ws.Range("H52").Value = -17.23313855436111
$ws.Range("G53").Value = -0.1397260522627475
$ws.Range("H53").Value = -10.84574677630037
$ws.Range("G54").Value = 0.1117980072647329
$ws.Range("H54").Value = 19.28633175351848
$ws.Range("G55").Value = 0.1107944668064794
$ws.Range("H55").Value = -2.021811762485136
$ws.Range("G56").Value = -0.01972015669265693
$ws.Range("H56").Value = -170.1105859959596
$ws.Range("G57").Value = -0.02750377474479362
$ws.Range("H57").Value = -20.27508715225616
$ws.Range("G58").Value = 0.04110046564293043
$ws.Range("H58").Value = -27.11129027174636
$ws.Range("G59").Value = 0.08097584729081388
$ws.Range("H59").Value = 12.74844160164536
$ws.Range("G60").Value = 0.06960502182728184
$ws.Range("H60").Value = -0.5272413504702349
$ws.Range("G61").Value = 0.08130271126037301
$ws.Range("H61").Value = 71.06713688578374
$ws.Range("G62").Value = 0.05208887078238711
$ws.Range("H62").Value = -28.60658311335892
$ws.Range("G63").Value = 0.0547461991295388
$ws.Range("H63").Value = -16.28190060695477
$ws.Range("G64").Value = -0.005402420019937308
$ws.Range("H64").Value = 86.95556244822697
$ws.Range("G65").Value = -0.003915631975341612
$ws.Range("H65").Value = 92.06135921508634
$ws.Range("G66").Value = 0.06823246015193106
$ws.Range("H66").Value = 260.3704910560059
$ws.Range("G67").Value = 0.05132886627199526
$ws.Range("H67").Value = 96.28715229445474
$ws.Range("G68").Value = -0.03571448348794981
$ws.Range("H68").Value = -6366.211408857946
$ws.Range("G69").Value = -0.01942987367203679
$ws.Range("H69").Value = -50.47561976556766
$ws.Range("G70").Value = -0.01992507918868098
$ws.Range("H70").Value = 27.40808167650369
$ws.Range("G71").Value = -0.06659038507724785
$ws.Range("H71").Value = -20.87135031682649
$ws.Range("G72").Value = -0.1377918717084664
$ws.Range("H72").Value = 7.107112381721993
$ws.Range("G73").Value = -0.1766709723679629
$ws.Range("H73").Value = -22.01890494863352
$ws.Range("G74").Value = 0.1467574933048519
$ws.Range("H74").Value = 16.42697822012114
$ws.Range("G75").Value = 0.1383480723548173
$ws.Range("H75").Value = 2.349669809232495
$ws.Range("G76").Value = -0.0878560261411637
$ws.Range("H76").Value = -155.1140204107573
$ws.Range("G77").Value = -0.0620567886087444
$ws.Range("H77").Value = -34.35983734667369
$ws.Range("G78").Value = 0.08608894621073446
$ws.Range("H78").Value = -6.599208482658879
$ws.Range("G79").Value = 0.1063874414704065
$ws.Range("H79").Value = 10.24099313551742
$ws.Range("G80").Value = -0.2137798715081172
$ws.Range("H80").Value = -31.62743001023166
$ws.Range("G81").Value = -0.2049090135341464
$ws.Range("H81").Value = 5.33126666261559
$ws.Range("G82").Value = 0.1652591466990643
$ws.Range("H82").Value = 19.1088855047139
$ws.Range("G83").Value = 0.1769682345928579
$ws.Range("H83").Value = 7.501318184906697
$ws.Range("G84").Value = 0.0191571265030936
$ws.Range("H84").Value = 36.88523599448757
$ws.Range("G85").Value = 0.04087498526755463
$ws.Range("H85").Value = 80.54973095015528
